$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20 (2025Q2): update total_customers and new_customers
$ws.Range("C20").Value = 349
$ws.Range("E20").Value = 79

# Row 21 (2025Q3): update total_customers, returning_customers, new_customers, recurrence_rate
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 9.169054441260744
